$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rowUpdates = @(
    @{ Row = 2; E = 3; F = 1; G = 4.873099; H = 14.619297; I = 0.01719598526069697; J = 0.01727593400119405; K = 3; L = 1; M = 1.009860666666667; N = 3.029582; O = 0.01353413605720072; P = 0.01542521070970148; Q = 4.921151004872667; R = 44.290359043854; S = 0.0002327328041558911; T = 0.0002664849221753144 },
    @{ Row = 3; E = 3; F = 1; G = 4.873099; H = 14.619297; I = 0.01719598526069697; J = 0.01727593400119405; K = 3; L = 1; M = 46.15376066666666; N = 138.461282; O = 0.6185519418990597; P = 0.704979911415303; Q = 224.9118449509726; R = 2024.206604558754; S = 0.01063661007587172; T = 0.0121791864217784 },
    @{ Row = 4; D = "Inflammatory-Mac"; E = 3; F = 1; G = 4.873099; H = 14.619297; I = 0.01719598526069697; J = 0.01727593400119405; K = 1; L = 0.3333333333333333; M = 0.009315666666666667; N = 0.027947; O = 0.0001248484115599408; P = 0.000142293017222847; Q = 0.04539616591766667; R = 0.408565493259; S = 0.000002146891445006172; T = 0.000002458244774372673 },
    @{ Row = 5; D = "MuSCs"; E = 3; F = 1; G = 4.873099; H = 14.619297; I = 0.01719598526069697; J = 0.01727593400119405; K = 2; L = 1; M = 27.4428835; N = 54.885767; O = 0.3677890736321797; P = 0.2794525848577725; Q = 133.7318881409665; R = 802.391328845799; S = 0.006324495489224356; T = 0.004827804412465958 },
    @{ Row = 6; E = 3; F = 1; G = 271.3121946666667; H = 813.936584; I = 0.9573949760789487; J = 0.9618461617095089; K = 3; L = 1; M = 1.009860666666667; N = 3.029582; O = 0.01353413605720072; P = 0.01542521070970148; Q = 273.9875137808765; R = 2465.887624027888; S = 0.01295751386673292; T = 0.01483667971468678 },
    @{ Row = 7; E = 3; F = 1; G = 271.3121946666667; H = 813.936584; I = 0.9573949760789487; J = 0.9618461617095089; K = 3; L = 1; M = 46.15376066666666; N = 138.461282; O = 0.6185519418990597; P = 0.704979911415303; Q = 12522.07809859341; R = 112698.7028873407; S = 0.5921985216180375; T = 0.6780822218771188 },
    @{ Row = 8; D = "Inflammatory-Mac"; E = 3; F = 1; G = 271.3121946666667; H = 813.936584; I = 0.9573949760789487; J = 0.9618461617095089; K = 1; L = 0.3333333333333333; M = 0.009315666666666667; N = 0.027947; O = 0.0001248484115599408; P = 0.000142293017222847; Q = 2.527453968116445; R = 22.747085713048; S = 0.0001195292419989243; T = 0.0001368639924538604 },
    @{ Row = 9; D = "MuSCs"; E = 3; F = 1; G = 271.3121946666667; H = 813.936584; I = 0.9573949760789487; J = 0.9618461617095089; K = 2; L = 1; M = 27.4428835; N = 54.885767; O = 0.3677890736321797; P = 0.2794525848577725; Q = 7445.588950366656; R = 44673.53370219993; S = 0.3521194113521794; T = 0.2687903961252493 },
    @{ Row = 10; E = 3; F = 1; G = 1.022486666666667; H = 3.067460000000001; I = 0.003608107622943672; J = 0.003624882681520371; K = 3; L = 1; M = 1.009860666666667; N = 3.029582; O = 0.01353413605720072; P = 0.01542521070970148; Q = 1.032569066857778; R = 9.293121601720001; S = 0.00004883261947794273; T = 0.00005591457916039944 },
    @{ Row = 11; E = 3; F = 1; G = 1.022486666666667; H = 3.067460000000001; I = 0.003608107622943672; J = 0.003624882681520371; K = 3; L = 1; M = 46.15376066666666; N = 138.461282; O = 0.6185519418990597; P = 0.704979911415303; Q = 47.19160489819111; R = 424.72444408372; S = 0.002231801976752608; T = 0.002555469471709097 },
    @{ Row = 12; D = "Inflammatory-Mac"; E = 3; F = 1; G = 1.022486666666667; H = 3.067460000000001; I = 0.003608107622943672; J = 0.003624882681520371; K = 1; L = 0.3333333333333333; M = 0.009315666666666667; N = 0.027947; O = 0.0001248484115599408; P = 0.000142293017222847; Q = 0.009525144957777778; R = 0.08572630462000001; S = 0.0000004504665054618312; T = 0.000000515795493832378 },
    @{ Row = 13; D = "MuSCs"; E = 3; F = 1; G = 1.022486666666667; H = 3.067460000000001; I = 0.003608107622943672; J = 0.003624882681520371; K = 2; L = 1; M = 27.4428835; N = 54.885767; O = 0.3677890736321797; P = 0.2794525848577725; Q = 28.05998247363667; R = 168.35989484182; S = 0.001327022560207659; T = 0.001012982835157041 },
    @{ Row = 14; E = 2; F = 1; G = 3.9343185; H = 7.868637; I = 0.01388325641175922; J = 0.009298535592467514; K = 3; L = 1; M = 1.009860666666667; N = 3.029582; O = 0.01353413605720072; P = 0.01542521070970148; Q = 3.973113503289; R = 23.838681019734; S = 0.0001878978811937536; T = 0.0001434318708054703 },
    @{ Row = 15; E = 2; F = 1; G = 3.9343185; H = 7.868637; I = 0.01388325641175922; J = 0.009298535592467514; K = 3; L = 1; M = 46.15376066666666; N = 138.461282; O = 0.6185519418990597; P = 0.704979911415303; Q = 181.583594435439; R = 1089.501566612634; S = 0.008587515213376237; T = 0.006555280798269791 },
    @{ Row = 16; D = "Inflammatory-Mac"; E = 2; F = 1; G = 3.9343185; H = 7.868637; I = 0.01388325641175922; J = 0.009298535592467514; K = 1; L = 0.3333333333333333; M = 0.009315666666666667; N = 0.027947; O = 0.0001248484115599408; P = 0.000142293017222847; Q = 0.0366507997065; R = 0.219904798239; S = 0.000001733302510287502; T = 0.000001323116685206236 },
    @{ Row = 17; D = "MuSCs"; E = 2; F = 1; G = 3.9343185; H = 7.868637; I = 0.01388325641175922; J = 0.009298535592467514; K = 2; L = 1; M = 27.4428835; N = 54.885767; O = 0.3677890736321797; P = 0.2794525848577725; Q = 107.9690442473947; R = 431.876176989579; S = 0.005106110014678942; T = 0.002598499806707046 },
    @{ Row = 18; E = 3; F = 1; G = 2.243757; H = 6.731271; I = 0.007917674625651083; J = 0.007954486015309181; K = 3; L = 1; M = 1.009860666666667; N = 3.029582; O = 0.01353413605720072; P = 0.01542521070970148; Q = 2.265881939858; R = 20.392937458722; S = 0.0001071588856402075; T = 0.0001226996228735178 },
    @{ Row = 19; E = 3; F = 1; G = 2.243757; H = 6.731271; I = 0.007917674625651083; J = 0.007954486015309181; K = 3; L = 1; M = 46.15376066666666; N = 138.461282; O = 0.6185519418990597; P = 0.704979911415303; Q = 103.557823572158; R = 932.0204121494219; S = 0.004897493015021388; T = 0.005607752846426933 },
    @{ Row = 20; D = "Inflammatory-Mac"; E = 3; F = 1; G = 2.243757; H = 6.731271; I = 0.007917674625651083; J = 0.007954486015309181; K = 1; L = 0.3333333333333333; M = 0.009315666666666667; N = 0.027947; O = 0.0001248484115599408; P = 0.000142293017222847; Q = 0.020902092293; R = 0.188118830637; S = 0.0000009885091002609866; T = 0.000001131867815575285 },
    @{ Row = 21; D = "MuSCs"; E = 3; F = 1; G = 2.243757; H = 6.731271; I = 0.007917674625651083; J = 0.007954486015309181; K = 2; L = 1; M = 27.4428835; N = 54.885767; O = 0.3677890736321797; P = 0.2794525848577725; Q = 61.5751619533095; R = 369.450971719857; S = 0.002912034215889227; T = 0.002222901678193154 }
)

foreach ($u in $rowUpdates) {
    $r = $u.Row
    if ($u.ContainsKey("D")) {
        $ws.Cells.Item($r, 4).Value = $u.D
    }
    $ws.Cells.Item($r, 5).Value = $u.E
    $ws.Cells.Item($r, 6).Value = $u.F
    $ws.Cells.Item($r, 7).Value = $u.G
    $ws.Cells.Item($r, 8).Value = $u.H
    $ws.Cells.Item($r, 9).Value = $u.I
    $ws.Cells.Item($r, 10).Value = $u.J
    $ws.Cells.Item($r, 11).Value = $u.K
    $ws.Cells.Item($r, 12).Value = $u.L
    $ws.Cells.Item($r, 13).Value = $u.M
    $ws.Cells.Item($r, 14).Value = $u.N
    $ws.Cells.Item($r, 15).Value = $u.O
    $ws.Cells.Item($r, 16).Value = $u.P
    $ws.Cells.Item($r, 17).Value = $u.Q
    $ws.Cells.Item($r, 18).Value = $u.R
    $ws.Cells.Item($r, 19).Value = $u.S
    $ws.Cells.Item($r, 20).Value = $u.T
}
